$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Modify existing "Antigua BD" values ---
# VALOR MORA total (E11): 56940 -> 113880
$ws.Range("E11").Value = 113880

# Cant. Periodos (F13): 1 -> 2 (now two EC periods listed below)
$ws.Range("F13").Value = 2

# --- Insert a new row for the new EC period (2509), pushing the
# signature block (rows 21-22) down to rows 22-23 ---
$ws.Rows.Item(17).Insert()

# Duplicate the formatting of the existing data row (16) into the
# newly inserted row (17)
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))

# Populate the new data row with the new period (2509), keeping the
# same worker/value data as the existing row
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73137696"
$ws.Range("D17").Value = "EDWIN ALFONSO CARABALLO POSADA"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500
